# Update the "Price" column (D) values for specific rows in the crypto
# price table to reflect the latest scrape from the GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "278.23"
    3  = "22.89"
    4  = "6.359"
    5  = "0.06264"
    6  = "3.661"
    7  = "6.633"
    8  = "1.400"
    9  = "0.8305"
    10 = "0.01388"
    12 = "0.08455"
    13 = "0.03500"
    14 = "0.03219"
    15 = "4.095"
    16 = "0.09298"
    17 = "0.001648"
    18 = "0.04746"
    19 = "0.006266"
    20 = "0.005736"
    22 = "0.0001498"
    23 = "3.728"
    25 = "0.3327"
    26 = "0.1259"
    28 = "0.0002698"
    40 = "0.04748"
    41 = "0.007100"
    42 = "0.1170"
    43 = "0.003650"
    44 = "0.01229"
    45 = "0.00006082"
    46 = "0.0009879"
    48 = "0.7804"
    50 = "0.00001399"
    51 = "0.01239"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    # Prefix with an apostrophe so the numeric-looking text is stored as a
    # literal string (matching the original cell's text type) instead of
    # being coerced into a numeric value.
    $cell.Value = "'" + $updates[$row]
    # Writing a quote-prefixed value marks the cell with a quote-prefix
    # style; reset back to Normal so the cell formatting matches the rest
    # of the untouched data cells in the column.
    $cell.Style = "Normal"
}
